$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 188.25
$ws.Range("I9").Value = 188.25
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 188.25
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -19.25
$ws.Range("N9").ClearContents()
$ws.Range("H18").Value = 1450
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H41").Value = 2699.6
$ws.Range("I41").Value = 1999.6666
$ws.Range("J41").Value = 3749.5
$ws.Range("K41").Value = 1999.6666
$ws.Range("L41").Value = 3749.5
$ws.Range("M41").Value = -1559.6666
$ws.Range("N41").Value = -4629.5
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H116").Value = 4333.3335
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4333.3335
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4333.3335
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -11217.3335
$ws.Range("H132").Value = 13902.579
$ws.Range("I132").Value = 13246.8125
$ws.Range("K132").Value = 39740.4375
$ws.Range("M132").Value = -37210.4375
$ws.Range("H138").Value = 2565.611
$ws.Range("I138").Value = 1853.5
$ws.Range("J138").Value = 3135.3
$ws.Range("K138").Value = 5560.5
$ws.Range("L138").Value = 9405.900000000001
$ws.Range("M138").Value = -420.5
$ws.Range("N138").Value = -19685.9

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2386460.5
$ws.Range("I32").Value = 3542.7778
$ws.Range("J32").Value = 16683966
$ws.Range("K32").Value = 3542.7778
$ws.Range("L32").Value = 16683966
$ws.Range("M32").Value = -3255.7778
$ws.Range("N32").Value = -16684540
$ws.Range("H45").Value = 2910
$ws.Range("I45").Value = 2005.6
$ws.Range("K45").Value = 2005.6
$ws.Range("M45").Value = -1628.6
$ws.Range("H94").Value = 18000
$ws.Range("J94").Value = 18000
$ws.Range("L94").Value = 18000
$ws.Range("N94").Value = -19802
$ws.Range("H132").Value = 3117.2307
$ws.Range("I132").Value = 1802.8
$ws.Range("K132").Value = 5408.4
$ws.Range("M132").Value = -2878.4

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3547.9285
$ws.Range("I20").Value = 4206.6
$ws.Range("K20").Value = 4206.6
$ws.Range("M20").Value = -3959.6
$ws.Range("H80").Value = 170.3
$ws.Range("J80").Value = 188.78572
$ws.Range("L80").Value = 188.78572
$ws.Range("N80").Value = -2184.78572
$ws.Range("H83").Value = 170.3
$ws.Range("J83").Value = 188.78572
$ws.Range("L83").Value = 943.9286
$ws.Range("N83").Value = -10927.9286
$ws.Range("H86").Value = 3706.889
$ws.Range("I86").Value = 1536.1428
$ws.Range("K86").Value = 1536.1428
$ws.Range("M86").Value = -413.1428000000001
$ws.Range("H89").Value = 3706.889
$ws.Range("I89").Value = 1536.1428
$ws.Range("K89").Value = 7680.714
$ws.Range("M89").Value = -2064.714

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 144.81818
$ws.Range("I7").Value = 121.28571
$ws.Range("J7").Value = 186
$ws.Range("K7").Value = 121.28571
$ws.Range("L7").Value = 186
$ws.Range("M7").Value = -8.285709999999995
$ws.Range("N7").Value = -412
$ws.Range("H22").Value = 1509.7368
$ws.Range("I22").Value = 898.63635
$ws.Range("J22").Value = 2350
$ws.Range("K22").Value = 898.63635
$ws.Range("L22").Value = 2350
$ws.Range("M22").Value = -548.63635
$ws.Range("N22").Value = -3050
$ws.Range("H31").Value = 8591.056
$ws.Range("I31").Value = 3298
$ws.Range("J31").Value = 9649.666999999999
$ws.Range("K31").Value = 3298
$ws.Range("L31").Value = 9649.666999999999
$ws.Range("M31").Value = -3003
$ws.Range("N31").Value = -10239.667
$ws.Range("H34").Value = 8591.056
$ws.Range("I34").Value = 3298
$ws.Range("J34").Value = 9649.666999999999
$ws.Range("K34").Value = 3298
$ws.Range("L34").Value = 9649.666999999999
$ws.Range("M34").Value = -3096
$ws.Range("N34").Value = -10053.667
$ws.Range("H58").Value = 5551.222
$ws.Range("I58").Value = 999
$ws.Range("J58").Value = 6851.857
$ws.Range("K58").Value = 999
$ws.Range("L58").Value = 6851.857
$ws.Range("M58").Value = -796
$ws.Range("N58").Value = -7257.857
$ws.Range("H59").Value = 33842
$ws.Range("I59").Value = 8500
$ws.Range("J59").Value = 41082.57
$ws.Range("K59").Value = 8500
$ws.Range("L59").Value = 41082.57
$ws.Range("M59").Value = -7355
$ws.Range("N59").Value = -43372.57
$ws.Range("H107").Value = 278.2353
$ws.Range("I107").Value = 285.14285
$ws.Range("K107").Value = 285.14285
$ws.Range("M107").Value = 1634.85715
$ws.Range("H136").Value = 5551.222
$ws.Range("I136").Value = 999
$ws.Range("J136").Value = 6851.857
$ws.Range("K136").Value = 2997
$ws.Range("L136").Value = 20555.571
$ws.Range("M136").Value = -447
$ws.Range("N136").Value = -25655.571

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.81818
$ws.Range("I2").Value = 34.210526
$ws.Range("J2").Value = 31.333334
$ws.Range("K2").Value = 205.263156
$ws.Range("L2").Value = 188.000004
$ws.Range("M2").Value = -92.26315600000001
$ws.Range("N2").Value = -414.000004
$ws.Range("H38").Value = 323.95
$ws.Range("I38").Value = 357.16666
$ws.Range("K38").Value = 1071.49998
$ws.Range("M38").Value = -724.4999800000001
$ws.Range("H131").Value = 1477.4445
$ws.Range("I131").Value = 1435.125
$ws.Range("J131").Value = 1816
$ws.Range("K131").Value = 4305.375
$ws.Range("L131").Value = 5448
$ws.Range("M131").Value = 734.625
$ws.Range("N131").Value = -15528
$ws.Range("H140").Value = 3979.8
$ws.Range("J140").Value = 5250
$ws.Range("L140").Value = 15750
$ws.Range("N140").Value = -26110

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11802.4
$ws.Range("I70").Value = 5503.5
$ws.Range("K70").Value = 5503.5
$ws.Range("M70").Value = -5233.5
$ws.Range("H73").Value = 11802.4
$ws.Range("I73").Value = 5503.5
$ws.Range("K73").Value = 5503.5
$ws.Range("M73").Value = -4567.5
$ws.Range("H80").Value = 2017.75
$ws.Range("I80").Value = 2057
$ws.Range("J80").Value = 1900
$ws.Range("K80").Value = 2057
$ws.Range("L80").Value = 1900
$ws.Range("M80").Value = -1059
$ws.Range("N80").Value = -3896
$ws.Range("H83").Value = 2017.75
$ws.Range("I83").Value = 2057
$ws.Range("J83").Value = 1900
$ws.Range("K83").Value = 10285
$ws.Range("L83").Value = 9500
$ws.Range("M83").Value = -5293
$ws.Range("N83").Value = -19484
$ws.Range("H93").Value = 45000
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 4476.5366
$ws.Range("I132").Value = 3957.7742
$ws.Range("K132").Value = 11873.3226
$ws.Range("M132").Value = -9343.3226

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3156.8
$ws.Range("I61").Value = 1395.6364
$ws.Range("K61").Value = 1395.6364
$ws.Range("M61").Value = -1193.6364
$ws.Range("H68").Value = 7977
$ws.Range("I68").Value = 6414.2
$ws.Range("J68").Value = 8535.143
$ws.Range("K68").Value = 6414.2
$ws.Range("L68").Value = 8535.143
$ws.Range("M68").Value = -5665.2
$ws.Range("N68").Value = -10033.143
$ws.Range("H71").Value = 7977
$ws.Range("I71").Value = 6414.2
$ws.Range("J71").Value = 8535.143
$ws.Range("K71").Value = 32071
$ws.Range("L71").Value = 42675.715
$ws.Range("M71").Value = -28327
$ws.Range("N71").Value = -50163.715
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H82").Value = 4090.3845
$ws.Range("I82").Value = 2245.5
$ws.Range("J82").Value = 5671.7144
$ws.Range("K82").Value = 2245.5
$ws.Range("L82").Value = 5671.7144
$ws.Range("M82").Value = -1884.5
$ws.Range("N82").Value = -6393.7144
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H85").Value = 4090.3845
$ws.Range("I85").Value = 2245.5
$ws.Range("J85").Value = 5671.7144
$ws.Range("K85").Value = 2245.5
$ws.Range("L85").Value = 5671.7144
$ws.Range("M85").Value = -997.5
$ws.Range("N85").Value = -8167.7144
$ws.Range("H113").Value = 3156.8
$ws.Range("I113").Value = 1395.6364
$ws.Range("K113").Value = 1395.6364
$ws.Range("M113").Value = 774.3635999999999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 71201.92
$ws.Range("I2").Value = 92500
$ws.Range("K2").Value = 92500
$ws.Range("M2").Value = -92388
$ws.Range("H4").Value = 92445.45
$ws.Range("I4").Value = 112950
$ws.Range("K4").Value = 112950
$ws.Range("M4").Value = -112837
$ws.Range("H62").Value = 12000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 12000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -13248
$ws.Range("H65").Value = 12000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 60000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -66240
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
